# Update the confusion-matrix cell text/values with the corrected numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 / Row 4 confusion-matrix values (shared strings get rebuilt in this
# set-order: C3, C4, D3, D4 -> matches target sharedStrings.xml ordering).
$ws.Range("C3").Value = "17 (0.8095)"
$ws.Range("C4").Value = "2 (0.1429)"
$ws.Range("D3").Value = "4 (0.1905)"
$ws.Range("D4").Value = "12 (0.8571)"

# Restyle: the header cells (C1:D1) pick up vertical centering, while the
# grid body (B2:D4) keeps horizontal centering only. The net rendered look
# is unchanged; it mirrors how the cell-style indices were swapped.
$ws.Range("C1:D1").VerticalAlignment = -4108  # xlCenter
$ws.Range("B2:D4").VerticalAlignment = -4107  # xlBottom (Excel default -> no explicit vertical attr)

# Move the active selection to G7, matching the saved selection in the file.
$ws.Range("G7").Select()
